$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new Adafruit IO reading as row 59 (same layout as the header
# row: Timestamp, Feed Key, Value, Latitude, Longitude, Elevation).
# The source data stores every column as text (even the numeric-looking
# "Value"/lat/long columns), so values that Excel would otherwise
# auto-coerce to a number/date are written with a leading apostrophe to
# force a literal text entry.
$ws.Range("A59").Value = "'2024-09-25T18:06:40Z"
$ws.Range("B59").Value = "temperature"
$ws.Range("C59").Value = "'25"
$ws.Range("D59").Value = "N/A"
$ws.Range("E59").Value = "N/A"
$ws.Range("F59").Value = "N/A"

# Drop the "quote prefix" formatting that Excel applies to force-texted
# cells so the new row doesn't pick up a style that differs from the rest
# of the sheet (which uses the default/general style throughout).
$ws.Range("A59:F59").ClearFormats()
